$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-12"

# Update the July row label
$ws.Range("A8").Value = "July (through 07-12)"

# Update the 2022 column (I) for April (row 5)
$ws.Range("I5").Value = 115

# Update July row (row 8) values for each year column B..I
$ws.Range("B8").Value = 16
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 22
$ws.Range("E8").Value = 29
$ws.Range("F8").Value = 17
$ws.Range("G8").Value = 43
$ws.Range("H8").Value = 57
$ws.Range("I8").Value = 68

# Update Total row (row 9) values for each year column B..I
$ws.Range("B9").Value = 141
$ws.Range("C9").Value = 271
$ws.Range("D9").Value = 412
$ws.Range("E9").Value = 382
$ws.Range("F9").Value = 268
$ws.Range("G9").Value = 515
$ws.Range("H9").Value = 817
$ws.Range("I9").Value = 874
